# Adds two new columns, I ("I0") and J ("IF"), to the first worksheet,
# populating header cells in row 1 and data values for rows 2-63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (style copied from the other header cells, e.g. column H)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-63: columnI, columnJ values
$data = @(
    @(2,6,7),
    @(3,7,7),
    @(4,7,7),
    @(5,9,9),
    @(6,9,9),
    @(7,6,6),
    @(8,8,8),
    @(9,9,9),
    @(10,8,8),
    @(11,8,8),
    @(12,7,7),
    @(13,5,6),
    @(14,6,6),
    @(15,7,7),
    @(16,7,7),
    @(17,6,6),
    @(18,8,8),
    @(19,8,9),
    @(20,6,6),
    @(21,8,8),
    @(22,9,9),
    @(23,8,8),
    @(24,7,7),
    @(25,8,8),
    @(26,8,8),
    @(27,9,9),
    @(28,7,7),
    @(29,8,8),
    @(30,8,8),
    @(31,8,8),
    @(32,6,6),
    @(33,9,9),
    @(34,8,8),
    @(35,8,8),
    @(36,7,8),
    @(37,7,7),
    @(38,7,8),
    @(39,9,9),
    @(40,11,11),
    @(41,8,8),
    @(42,9,9),
    @(43,8,8),
    @(44,9,9),
    @(45,8,9),
    @(46,7,8),
    @(47,9,9),
    @(48,9,9),
    @(49,7,7),
    @(50,8,8),
    @(51,8,8),
    @(52,8,8),
    @(53,7,7),
    @(54,8,8),
    @(55,7,7),
    @(56,9,9),
    @(57,7,7),
    @(58,8,8),
    @(59,7,7),
    @(60,7,7),
    @(61,5,5),
    @(62,7,7),
    @(63,5,5)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
